# digital_economic_outlook_2024.xlsx
# Replace the "Country Code" column (ISO3 codes) with an "OECD Membership"
# (Yes/No) column on the "NDS development_e" sheet, and rename the
# "Country Name" header to "Country".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row -------------------------------------------------------
# Write the new B1 header first so that the shared-string table is built
# in the same order as the target workbook (OECD Membership, No, Yes, Country).
$ws.Cells.Item(1, 2).Value = "OECD Membership"

# --- Column B data: OECD membership Yes/No ----------------------------
# Countries that are NOT OECD members -> "No" (keep their original,
# non-wrapped cell style).
$noRows = @(2, 13, 16, 19)
foreach ($r in $noRows) {
    $ws.Cells.Item($r, 2).Value = "No"
}

# Countries that ARE OECD members -> "Yes" (wrap the text in the cell,
# matching the author's formatting of the new column).
$yesRows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 14, 15, 17, 18, 20, 21, 22, 23)
foreach ($r in $yesRows) {
    $ws.Cells.Item($r, 2).Value = "Yes"
    $ws.Cells.Item($r, 2).WrapText = $true
    $ws.Cells.Item($r, 2).VerticalAlignment = -4107
}

# A1 header renamed last so it becomes the final new shared string.
$ws.Cells.Item(1, 1).Value = "Country"

# --- Column widths ------------------------------------------------------
# Column B becomes wide enough to fit "OECD Membership" (stored width 15);
# column C is left untouched so it keeps its original (shared) width.
$ws.Columns.Item(2).ColumnWidth = 14.1

# --- Selection ------------------------------------------------------
# Reset the stored selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
